# Generate Report for Handback
# Refresh the handoff/handback/generate timestamps for the
# "bda00650-785f-46ab-84fa-0e9e1ebdd06c" row (row 2) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 11:14:41"

# --- zh-cn sheet ---
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 11:14:36"
$wsZhCn.Range("K2").Value = "2016-09-01 11:14:54"

# --- de-de sheet ---
# "Correspond Handoff Datetime" mirrors the Overview generate date,
# "Correspond Handback DateTime" is its own value.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 11:14:41"
$wsDeDe.Range("K2").Value = "2016-09-01 11:15:03"
